# Insert a new data row at sheet row 252 (pushes existing rows 252-314 down
# to 253-315) and populate it with the new weekly Albahaca price record for
# "Vega Modelo de Temuco".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("252:252").Insert()

$ws.Range("A252").Value = 10
$ws.Range("B252").Value = 'Vega Modelo de Temuco'
$ws.Range("C252").Value = 'La Araucanía'
$ws.Range("D252").Value = 44943
$ws.Range("E252").Value = 9
$ws.Range("F252").Value = 100112052
$ws.Range("G252").Value = 'Albahaca'
$ws.Range("H252").Value = 'Sin especificar'
$ws.Range("I252").Value = 'Primera'
$ws.Range("J252").Value = 35
$ws.Range("K252").Value = 5000
$ws.Range("L252").Value = 5000
$ws.Range("M252").Value = 5000
$ws.Range("N252").Value = '$/paquete'
$ws.Range("O252").Value = 'Región de La Araucanía'
$ws.Range("P252").Value = 5000
$ws.Range("Q252").Value = 1
$ws.Range("R252").Value = 'Hortaliza'
